$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: Category | PlayersNumber | Round Name | Points ---
$ws.Cells.Item(1,1).Value = "Category"
$ws.Cells.Item(1,2).Value = "PlayersNumber"
$ws.Cells.Item(1,3).Value = "Round Name"
$ws.Cells.Item(1,4).Value = "Points"

# Right-align the Round Name header (matches new bold+right style)
$ws.Cells.Item(1,3).HorizontalAlignment = -4152

# --- Data rows: round names remapped to short codes, values corrected, ---
# --- rows reordered into the canonical round sequence per category/draw ---
$data = @(
  @(2, "ATP-250", 32, "Q1", 0),
  @(3, "ATP-250", 32, "Q2", 0),
  @(4, "ATP-250", 32, "Q3", 10),
  @(5, "ATP-250", 32, "Q", 0),
  @(6, "ATP-250", 32, "R128", 0),
  @(7, "ATP-250", 32, "R64", 0),
  @(8, "ATP-250", 32, "R32", 20),
  @(9, "ATP-250", 32, "R16", 45),
  @(10, "ATP-250", 32, "QF", 90),
  @(11, "ATP-250", 32, "SF", 150),
  @(12, "ATP-250", 32, "F", 250),
  @(13, "ATP-250", 48, "Q1", 6),
  @(14, "ATP-250", 48, "Q2", 0),
  @(15, "ATP-250", 48, "Q3", 20),
  @(16, "ATP-250", 48, "Q", 0),
  @(17, "ATP-250", 48, "R128", 0),
  @(18, "ATP-250", 48, "R64", 10),
  @(19, "ATP-250", 48, "R32", 20),
  @(20, "ATP-250", 48, "R16", 45),
  @(21, "ATP-250", 48, "QF", 90),
  @(22, "ATP-250", 48, "SF", 150),
  @(23, "ATP-250", 48, "F", 250),
  @(24, "ATP-500", 32, "Q1", 3),
  @(25, "ATP-500", 32, "Q2", 0),
  @(26, "ATP-500", 32, "Q3", 10),
  @(27, "ATP-500", 32, "Q", 0),
  @(28, "ATP-500", 32, "R128", 0),
  @(29, "ATP-500", 32, "R64", 0),
  @(30, "ATP-500", 32, "R32", 45),
  @(31, "ATP-500", 32, "R16", 90),
  @(32, "ATP-500", 32, "QF", 180),
  @(33, "ATP-500", 32, "SF", 300),
  @(34, "ATP-500", 32, "F", 500),
  @(35, "ATP-500", 48, "Q1", 10),
  @(36, "ATP-500", 48, "Q2", 0),
  @(37, "ATP-500", 48, "Q3", 0),
  @(38, "ATP-500", 48, "Q", 0),
  @(39, "ATP-500", 48, "R128", 0),
  @(40, "ATP-500", 48, "R64", 20),
  @(41, "ATP-500", 48, "R32", 45),
  @(42, "ATP-500", 48, "R16", 90),
  @(43, "ATP-500", 48, "QF", 180),
  @(44, "ATP-500", 48, "SF", 300),
  @(45, "ATP-500", 48, "F", 500),
  @(46, "ATP-Fs", 8, "Q1", 200),
  @(47, "ATP-Fs", 8, "Q2", 200),
  @(48, "ATP-Fs", 8, "Q3", 200),
  @(49, "ATP-Fs", 8, "Q", 200),
  @(50, "ATP-Fs", 8, "R128", 200),
  @(51, "ATP-Fs", 8, "R64", 200),
  @(52, "ATP-Fs", 8, "R32", 200),
  @(53, "ATP-Fs", 8, "R16", 200),
  @(54, "ATP-Fs", 8, "QF", 600),
  @(55, "ATP-Fs", 8, "SF", 1000),
  @(56, "ATP-Fs", 8, "F", 1500),
  @(57, "Grand Slam", 128, "Q1", 8),
  @(58, "Grand Slam", 128, "Q2", 16),
  @(59, "Grand Slam", 128, "Q3", 25),
  @(60, "Grand Slam", 128, "Q", 10),
  @(61, "Grand Slam", 128, "R128", 45),
  @(62, "Grand Slam", 128, "R64", 90),
  @(63, "Grand Slam", 128, "R32", 180),
  @(64, "Grand Slam", 128, "R16", 360),
  @(65, "Grand Slam", 128, "QF", 720),
  @(66, "Grand Slam", 128, "SF", 1200),
  @(67, "Grand Slam", 128, "F", 2000),
  @(68, "Masters-1000", 56, "Q1", 4),
  @(69, "Masters-1000", 56, "Q2", 0),
  @(70, "Masters-1000", 56, "Q3", 25),
  @(71, "Masters-1000", 56, "Q", 0),
  @(72, "Masters-1000", 56, "R128", 10),
  @(73, "Masters-1000", 56, "R64", 45),
  @(74, "Masters-1000", 56, "R32", 90),
  @(75, "Masters-1000", 56, "R16", 180),
  @(76, "Masters-1000", 56, "QF", 360),
  @(77, "Masters-1000", 56, "SF", 600),
  @(78, "Masters-1000", 56, "F", 1000),
  @(79, "Masters-1000", 96, "Q1", 16),
  @(80, "Masters-1000", 96, "Q2", 0),
  @(81, "Masters-1000", 96, "Q3", 16),
  @(82, "Masters-1000", 96, "Q", 10),
  @(83, "Masters-1000", 96, "R128", 25),
  @(84, "Masters-1000", 96, "R64", 45),
  @(85, "Masters-1000", 96, "R32", 90),
  @(86, "Masters-1000", 96, "R16", 180),
  @(87, "Masters-1000", 96, "QF", 360),
  @(88, "Masters-1000", 96, "SF", 600),
  @(89, "Masters-1000", 96, "F", 1000)
)

foreach ($item in $data) {
  $r = $item[0]
  $ws.Cells.Item($r,1).Value = $item[1]
  $ws.Cells.Item($r,2).Value = $item[2]
  $ws.Cells.Item($r,3).Value = $item[3]
  $ws.Cells.Item($r,4).Value = $item[4]
  $ws.Cells.Item($r,3).HorizontalAlignment = -4152
}

# --- Column widths: add a width for column A, keep B/C as before ---
$ws.Columns.Item(1).ColumnWidth = 11.83

# --- Selection moved by the author while reviewing the data ---
$ws.Range("E68").Select()